# wzrv/expdata/2020.xlsx — "Add files via upload"
#
# The column "pt_max" (column F, every data row holds the constant 50) is
# removed from the table entirely, shifting every column to its right
# (boson, value, stat_u, syst_u, obs, diff) one position to the left.
# The shared string "pt_max" disappears from the workbook along with it.
#
# The header row (row 1) additionally becomes bold (on top of the
# pre-existing centered alignment), and the sheet's active selection moves
# to L19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the pt_max column (F) - Excel shifts G:L left to F:K automatically,
# re-numbering every cell reference, updating <dimension>, <cols>, shared
# strings, etc.
$ws.Range("F1").EntireColumn.Delete()

# Header row (now A1:K1) becomes bold, keeping its existing centered
# alignment.
$ws.Range("A1:K1").Font.Bold = $true

# Match the saved selection state.
$ws.Range("L19").Select() | Out-Null
